# Update vm_pu.xlsx results for 380 kV case (B2 setpoint 1.05 -> 1.02, with recomputed bus voltages)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041228053287797
$ws.Cells.Item(2, 4).Value = 1.054017009620546
$ws.Cells.Item(2, 5).Value = 1.049874139649618
$ws.Cells.Item(2, 6).Value = 1.061564441082372
$ws.Cells.Item(2, 9).Value = 1.045683326077794
$ws.Cells.Item(2, 10).Value = 1.046310394696135
$ws.Cells.Item(2, 11).Value = 1.056761380854792
$ws.Cells.Item(2, 12).Value = 1.052629984699511
$ws.Cells.Item(2, 13).Value = 1.064288158697128
$ws.Cells.Item(2, 14).Value = 1.047796275470083

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042115129792276
$ws.Cells.Item(3, 4).Value = 1.054651646024352
$ws.Cells.Item(3, 5).Value = 1.050643602167365
$ws.Cells.Item(3, 6).Value = 1.062371288255621
$ws.Cells.Item(3, 9).Value = 1.045903933372039
$ws.Cells.Item(3, 10).Value = 1.046843724536439
$ws.Cells.Item(3, 11).Value = 1.057209633772944
$ws.Cells.Item(3, 12).Value = 1.053211891331519
$ws.Cells.Item(3, 13).Value = 1.064909671141698
$ws.Cells.Item(3, 14).Value = 1.048330362699934

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042689724578381
$ws.Cells.Item(4, 4).Value = 1.055062603806165
$ws.Cells.Item(4, 5).Value = 1.051142374891397
$ws.Cells.Item(4, 6).Value = 1.062894223798939
$ws.Cells.Item(4, 9).Value = 1.046045557285307
$ws.Cells.Item(4, 10).Value = 1.04718876797049
$ws.Cells.Item(4, 11).Value = 1.057499298104323
$ws.Cells.Item(4, 12).Value = 1.053588639707767
$ws.Cells.Item(4, 13).Value = 1.065312025755028
$ws.Cells.Item(4, 14).Value = 1.048675896135239

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042931425686584
$ws.Cells.Item(5, 4).Value = 1.055235441806334
$ws.Cells.Item(5, 5).Value = 1.051352267632462
$ws.Cells.Item(5, 6).Value = 1.063114267847946
$ws.Cells.Item(5, 9).Value = 1.046104826133479
$ws.Cells.Item(5, 10).Value = 1.047333809447053
$ws.Cells.Item(5, 11).Value = 1.057620979411313
$ws.Cells.Item(5, 12).Value = 1.053747075028565
$ws.Cells.Item(5, 13).Value = 1.06548122063962
$ws.Cells.Item(5, 14).Value = 1.048821143587334

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042972016619228
$ws.Cells.Item(6, 4).Value = 1.055264466198301
$ws.Cells.Item(6, 5).Value = 1.051387521733859
$ws.Cells.Item(6, 6).Value = 1.063151225999077
$ws.Cells.Item(6, 9).Value = 1.046114761786365
$ws.Cells.Item(6, 10).Value = 1.047358161646701
$ws.Cells.Item(6, 11).Value = 1.057641404711713
$ws.Cells.Item(6, 12).Value = 1.053773679925841
$ws.Cells.Item(6, 13).Value = 1.065509631799011
$ws.Cells.Item(6, 14).Value = 1.048845530369899

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.042692953646394
$ws.Cells.Item(7, 4).Value = 1.055064912997912
$ws.Cells.Item(7, 5).Value = 1.051145178670678
$ws.Cells.Item(7, 6).Value = 1.062897163247166
$ws.Cells.Item(7, 9).Value = 1.046046350299882
$ws.Cells.Item(7, 10).Value = 1.047190706081422
$ws.Cells.Item(7, 11).Value = 1.057500924385285
$ws.Cells.Item(7, 12).Value = 1.05359075653235
$ws.Cells.Item(7, 13).Value = 1.065314286370477
$ws.Cells.Item(7, 14).Value = 1.048677836998511

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041527720462788
$ws.Cells.Item(8, 4).Value = 1.05423142360221
$ws.Cells.Item(8, 5).Value = 1.050134000049754
$ws.Cells.Item(8, 6).Value = 1.061836941419662
$ws.Cells.Item(8, 9).Value = 1.04575811353695
$ws.Cells.Item(8, 10).Value = 1.046490647003174
$ws.Cells.Item(8, 11).Value = 1.056912948996223
$ws.Cells.Item(8, 12).Value = 1.052826597244129
$ws.Cells.Item(8, 13).Value = 1.064498160527629
$ws.Cells.Item(8, 14).Value = 1.047976783756074

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.0394790672118
$ws.Cells.Item(9, 4).Value = 1.052765129199632
$ws.Cells.Item(9, 5).Value = 1.048358983796992
$ws.Cells.Item(9, 6).Value = 1.059975300159806
$ws.Cells.Item(9, 9).Value = 1.045241632931594
$ws.Cells.Item(9, 10).Value = 1.045256673818058
$ws.Cells.Item(9, 11).Value = 1.055873969675356
$ws.Cells.Item(9, 12).Value = 1.051481765143126
$ws.Cells.Item(9, 13).Value = 1.06306160105592
$ws.Cells.Item(9, 14).Value = 1.046741058187496

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038116501991059
$ws.Cells.Item(10, 4).Value = 1.051789334728377
$ws.Cells.Item(10, 5).Value = 1.047180322345313
$ws.Cells.Item(10, 6).Value = 1.058738759796873
$ws.Cells.Item(10, 9).Value = 1.044891598600996
$ws.Cells.Item(10, 10).Value = 1.044433842505318
$ws.Cells.Item(10, 11).Value = 1.055179452437771
$ws.Cells.Item(10, 12).Value = 1.050586444840766
$ws.Cells.Item(10, 13).Value = 1.062105035367841
$ws.Cells.Item(10, 14).Value = 1.045917058359912

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037527274079533
$ws.Cells.Item(11, 4).Value = 1.051367239245551
$ws.Cells.Item(11, 5).Value = 1.046671081577274
$ws.Cells.Item(11, 6).Value = 1.058204428080404
$ws.Cells.Item(11, 9).Value = 1.044738684767538
$ws.Cells.Item(11, 10).Value = 1.044077519432564
$ws.Cells.Item(11, 11).Value = 1.054878293398368
$ws.Cells.Item(11, 12).Value = 1.050199071282427
$ws.Cells.Item(11, 13).Value = 1.061691121933991
$ws.Cells.Item(11, 14).Value = 1.045560229267524

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037308525984114
$ws.Cells.Item(12, 4).Value = 1.051210520624138
$ws.Cells.Item(12, 5).Value = 1.046482098039618
$ws.Cells.Item(12, 6).Value = 1.058006120439089
$ws.Cells.Item(12, 9).Value = 1.044681684238565
$ws.Cells.Item(12, 10).Value = 1.043945161634722
$ws.Cells.Item(12, 11).Value = 1.054766366489955
$ws.Cells.Item(12, 12).Value = 1.050055230856532
$ws.Cells.Item(12, 13).Value = 1.061537420566036
$ws.Cells.Item(12, 14).Value = 1.045427683506427

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03735544285818
$ws.Cells.Item(13, 4).Value = 1.051244134264648
$ws.Cells.Item(13, 5).Value = 1.046522627876233
$ws.Cells.Item(13, 6).Value = 1.058048650515644
$ws.Cells.Item(13, 9).Value = 1.044693920152037
$ws.Cells.Item(13, 10).Value = 1.04397355297246
$ws.Cells.Item(13, 11).Value = 1.054790378022523
$ws.Cells.Item(13, 12).Value = 1.050086082950769
$ws.Cells.Item(13, 13).Value = 1.061570387991818
$ws.Cells.Item(13, 14).Value = 1.045456115163121

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037509189890804
$ws.Cells.Item(14, 4).Value = 1.051354283469555
$ws.Cells.Item(14, 5).Value = 1.046655456637273
$ws.Cells.Item(14, 6).Value = 1.058188032492076
$ws.Cells.Item(14, 9).Value = 1.044733977192419
$ws.Cells.Item(14, 10).Value = 1.044066578760949
$ws.Cells.Item(14, 11).Value = 1.054869042763271
$ws.Cells.Item(14, 12).Value = 1.0501871804138
$ws.Cells.Item(14, 13).Value = 1.061678416013888
$ws.Cells.Item(14, 14).Value = 1.045549273058901

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037603934109533
$ws.Cells.Item(15, 4).Value = 1.051422158888229
$ws.Cells.Item(15, 5).Value = 1.046737319553973
$ws.Cells.Item(15, 6).Value = 1.058273932507863
$ws.Cells.Item(15, 9).Value = 1.044758630977203
$ws.Cells.Item(15, 10).Value = 1.044123894586358
$ws.Cells.Item(15, 11).Value = 1.054917502394129
$ws.Cells.Item(15, 12).Value = 1.050249476214604
$ws.Cells.Item(15, 13).Value = 1.061744981588708
$ws.Cells.Item(15, 14).Value = 1.045606670279356

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038155623490229
$ws.Cells.Item(16, 4).Value = 1.051817357031198
$ws.Cells.Item(16, 5).Value = 1.04721414288296
$ws.Cells.Item(16, 6).Value = 1.058774244939851
$ws.Cells.Item(16, 9).Value = 1.044901718684012
$ws.Cells.Item(16, 10).Value = 1.044457489930442
$ws.Cells.Item(16, 11).Value = 1.055199430458288
$ws.Cells.Item(16, 12).Value = 1.050612160113982
$ws.Cells.Item(16, 13).Value = 1.062132511595175
$ws.Cells.Item(16, 14).Value = 1.045940739367092

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038501891309157
$ws.Cells.Item(17, 4).Value = 1.052065370796917
$ws.Cells.Item(17, 5).Value = 1.047513544655675
$ws.Cells.Item(17, 6).Value = 1.059088373147571
$ws.Cells.Item(17, 9).Value = 1.044991113846083
$ws.Cells.Item(17, 10).Value = 1.04466673785241
$ws.Cells.Item(17, 11).Value = 1.055376162714163
$ws.Cells.Item(17, 12).Value = 1.050839745053889
$ws.Cells.Item(17, 13).Value = 1.062375676560635
$ws.Cells.Item(17, 14).Value = 1.04615028444511

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038703937843012
$ws.Cells.Item(18, 4).Value = 1.052210074392505
$ws.Cells.Item(18, 5).Value = 1.047688289304769
$ws.Cells.Item(18, 6).Value = 1.05927170466924
$ws.Cells.Item(18, 9).Value = 1.045043126545121
$ws.Cells.Item(18, 10).Value = 1.044788785486377
$ws.Cells.Item(18, 11).Value = 1.055479206215574
$ws.Cells.Item(18, 12).Value = 1.050972520934732
$ws.Cells.Item(18, 13).Value = 1.062517538043401
$ws.Cells.Item(18, 14).Value = 1.046272505400716

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038772843061317
$ws.Cells.Item(19, 4).Value = 1.052259421538892
$ws.Cells.Item(19, 5).Value = 1.047747891090973
$ws.Cells.Item(19, 6).Value = 1.059334233876315
$ws.Cells.Item(19, 9).Value = 1.045060839459145
$ws.Cells.Item(19, 10).Value = 1.044830400006168
$ws.Cells.Item(19, 11).Value = 1.055514334317944
$ws.Cells.Item(19, 12).Value = 1.051017799029628
$ws.Cells.Item(19, 13).Value = 1.062565913755034
$ws.Cells.Item(19, 14).Value = 1.0463141790179

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03846473234611
$ws.Cells.Item(20, 4).Value = 1.052038756975491
$ws.Cells.Item(20, 5).Value = 1.047481410426883
$ws.Cells.Item(20, 6).Value = 1.059054659201028
$ws.Cells.Item(20, 9).Value = 1.044981536030582
$ws.Cells.Item(20, 10).Value = 1.044644287846549
$ws.Cells.Item(20, 11).Value = 1.055357205291628
$ws.Cells.Item(20, 12).Value = 1.050815324289012
$ws.Cells.Item(20, 13).Value = 1.062349584408196
$ws.Cells.Item(20, 14).Value = 1.046127802557665

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037463911969809
$ws.Cells.Item(21, 4).Value = 1.051321845425444
$ws.Cells.Item(21, 5).Value = 1.046616337138779
$ws.Cells.Item(21, 6).Value = 1.058146983349962
$ws.Cells.Item(21, 9).Value = 1.044722186949026
$ws.Cells.Item(21, 10).Value = 1.044039185072885
$ws.Cells.Item(21, 11).Value = 1.054845879686995
$ws.Cells.Item(21, 12).Value = 1.050157408412219
$ws.Cells.Item(21, 13).Value = 1.061646603212601
$ws.Cells.Item(21, 14).Value = 1.045521840468657

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036835336347949
$ws.Cells.Item(22, 4).Value = 1.050871479957535
$ws.Cells.Item(22, 5).Value = 1.046073422370452
$ws.Cells.Item(22, 6).Value = 1.057577258385075
$ws.Cells.Item(22, 9).Value = 1.044557958335087
$ws.Cells.Item(22, 10).Value = 1.043658713044686
$ws.Cells.Item(22, 11).Value = 1.05452402496513
$ws.Cells.Item(22, 12).Value = 1.049744025290675
$ws.Cells.Item(22, 13).Value = 1.061204869335924
$ws.Cells.Item(22, 14).Value = 1.045140828126538

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.037168491256012
$ws.Cells.Item(23, 4).Value = 1.051110190080724
$ws.Cells.Item(23, 5).Value = 1.046361137220061
$ws.Cells.Item(23, 6).Value = 1.057879188078607
$ws.Cells.Item(23, 9).Value = 1.044645129245003
$ws.Cells.Item(23, 10).Value = 1.04386040994524
$ws.Cells.Item(23, 11).Value = 1.0546946802993
$ws.Cells.Item(23, 12).Value = 1.049963141004784
$ws.Cells.Item(23, 13).Value = 1.061439015814011
$ws.Cells.Item(23, 14).Value = 1.045342811459825

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038481522652372
$ws.Cells.Item(24, 4).Value = 1.052050782488245
$ws.Cells.Item(24, 5).Value = 1.047495930166149
$ws.Cells.Item(24, 6).Value = 1.059069892755442
$ws.Cells.Item(24, 9).Value = 1.044985864234861
$ws.Cells.Item(24, 10).Value = 1.044654432047686
$ws.Cells.Item(24, 11).Value = 1.055365771463043
$ws.Cells.Item(24, 12).Value = 1.050826358890693
$ws.Cells.Item(24, 13).Value = 1.062361374244064
$ws.Cells.Item(24, 14).Value = 1.046137961164732

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040008134855839
$ws.Cells.Item(25, 4).Value = 1.053143903499695
$ws.Cells.Item(25, 5).Value = 1.048817050313768
$ws.Cells.Item(25, 6).Value = 1.060455784671124
$ws.Cells.Item(25, 9).Value = 1.045376166218609
$ws.Cells.Item(25, 10).Value = 1.04557572296159
$ws.Cells.Item(25, 11).Value = 1.056142905174845
$ws.Cells.Item(25, 12).Value = 1.051829224894127
$ws.Cells.Item(25, 13).Value = 1.063432791857172
$ws.Cells.Item(25, 14).Value = 1.047060560417407

Write-Host "Updated vm_pu values for 380 kV case"